# Rename the worksheets and switch the active tab from "Connections" (was
# Sheet1) to "Structure" (was Sheet2), matching the authored diff.

$wb = $excel.ActiveWorkbook

$wsConnections = $wb.Worksheets.Item("Sheet1")
$wsConnections.Name = "Connections"

$wsStructure = $wb.Worksheets.Item("Sheet2")
$wsStructure.Name = "Structure"

# Make "Structure" the active / selected sheet (tabSelected moves to it,
# workbook view activeTab becomes 1).
$wsStructure.Activate()
$wsStructure.Select()
